# Apply "Add data for 2022-06-04" update to the carjacking-by-neighborhood-by-month workbook.
#
# Changes:
#  - Rename the sheet / update the "through" date from May 26 to May 27, 2022
#    (both the sheet tab name and the column-B header text).
#  - Bump a handful of monthly neighborhood carjacking counts (column B = current
#    "May 2022" column reflects the extra day's data; other touched columns are
#    pre-existing historical cells that simply needed incrementing / newly populating).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet tab and update the shared header text -----------------------
$ws.Name = "Through 2022-05-27"
$ws.Range("B1").Value = "May 2022 (through May 27)"

# --- Update / add individual cell values ---------------------------------------
# Row 3 - Austin
$ws.Range("Q3").Value = 7

# Row 4 - Humboldt Park
$ws.Range("V4").Value = 5

# Row 5 - Garfield Park
$ws.Range("G5").Value = 5
$ws.Range("L5").Value = 6

# Row 7 - North Lawndale
$ws.Range("AF7").Value = 3

# Row 8 - South Shore
$ws.Range("Q8").Value = 4
$ws.Range("AA8").Value = 3

# Row 13 - Washington Heights (new value)
$ws.Range("L13").Value = 1

# Row 16 - South Chicago
$ws.Range("L16").Value = 2

# Row 17 - Pullman (new value)
$ws.Range("L17").Value = 1

# Row 28 - West Town (new value)
$ws.Range("L28").Value = 1

# Row 30 - West Loop (new value)
$ws.Range("G30").Value = 1

# Row 32 - United Center
$ws.Range("AA32").Value = 2

# Row 39 - New City
$ws.Range("V39").Value = 2

# Row 46 - Little Village
$ws.Range("B46").Value = 2

# Row 92 - West Elsdon (new value)
$ws.Range("G92").Value = 1
